$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update encrypted note content stored in column B for rows 14 and 15
$ws.Range("B14").Value2 = "U2FsdGVkX1/kl9dHvj1AAodq6aaqyy8BUPPOOUxkT9Jf4uH3m9zzXNlXQtHolSAD4+idhFaI+E0y2axksFZkvTpkbZvl1z+cWGbs1pgftbkPR9FEcTtOVuvFjnAy+NrH16NYHPs252hZ+q6TruM+kQjn0Vuh4UjrM59MQTLPfEtZcFKsEYs176y0GH60M6Z9KZwPaYstUtXQsLEdIbI/Esh6BDqbkTtjw6qMv7f1+6XSJu0OhKtFJ5q9OrappOQqwkN2VcD/iy25+MdkmbV6Mbr8T8dPGqcDnk6U6Wzp1WcC5DcPzbXn6V0r0F9wH5IDo16z+DvuLDzFpkbOVRkt2QNOg40Oi5PIKMuyll8flS+Q+5zyNUn2bXzhnnjrQcRNhFrAmWUDhp+p/xFlExoOdJnWmjapeg1X21pw3ru+KCk1YqOd0MhXPr8e9y7SgMUWPGafVFyjxQ1QYEPmHkM++6puxhFzUGbAh+sh9Gs5aYYI8edJEpgNSE5TtnpUVr08F9DVCpl6CnwWgtU7CU8GxuWUblP37doiRkhLgBkIh+zCF3w0EVScZ6U779UZ10lSYMJlgMFDkFrnsEhE91FoIefVeCOuXliBb9zGETb+Eqr8kNViP5Ump4vJf2Iq8ZawPIDgKVS2X1mUfqMjHg1K50NdWrhxyTNjQTQvu3v0YIHDEsisWGJZcAsfMYb6tBSJTHSs2IFFlNh6X1YCq/HTkXiSh0vFZf8fzqW6kDkZ7Bsk12l2l1zp1r1fpZZF7Y+/iVH8wPmiI+R+7W0ecIFZH+7RVQcSO94UR7qeqmOlLVZOitDIhWk2U9U2+03puYdJ//xCWrazTSr04GZ8k3UO7KmLZ2MiVZMbNpZN3oRfaBmNwz4rvtRtlVgYyEOPla8XbbKVDLufxfH9+lNjj/HKDszGrReznh1B8sgBge7QVN/RfXMGWbQ+6SvkdLcn3rIo9VR7fSnPhnwA0YFRiuQWWxiUsJMjVeqrkh65z9f3JLYgXkTAABy5fIMM5VuzR8r4VTdFBYY581b2KtcChMaPS4G1xMA3y5WZ7RUin9FWo602xB6OMpCYFcmjz10RJyyIO4i+cJeLFwlyLxHnY9KtSmm/Psfs7CvzZ2XS6J7fHWPPCIn/XyIdNaWRsKs36WYlwFLiYTV0SLvhjUYLQca11FjMUnpdkuQ0X1mYXeYMJJd45kwA2XxAqBTxFttZErJB7JVjwyqeoy7co/tCGeHULTfZzlv1AQEMPKXWLJM9dSDeamf9lb1ScfHyQ2W3haGDK4bwTKndAsLAJOXkVNXJ7BNZPCQnQcQJUBmw4W7BdhQi2HrC25avYpFo45LUymByBKOUp+dv+cEIoz4O7afUZ0VD/DomYN9ZLGaF+NcVsBHmFq+OTRl8FgSG/V+b37L/bPNEef/wGMltliKjs2DOyQ=="
$ws.Range("B15").Value2 = "U2FsdGVkX19Ae8VO010rhTmO9aJs+DbkAtfYaJE5ZdPDtKn0twktiq0W0ZsjEdDuqK2N4+W1eCmgLR5xnuwbPlyUh9uc9s9hoeSQqF0yv/99x7wS3JXeOAwUgpvQzOXvjAw4e/48PccHObB9v4+YSFvk2dqqqETtUEtM+iWe45IICdtgxjaOg4BYoQSzsKNdFg7IPlIv/U2tAgi0NRgwBiXniL4aZqsW6nCEtV+A6QPlEpRa/wuJLRZ9wjB0rGKd"

# Move the active cell selection from B17 to B19 (frozen-pane bottom-left pane)
$ws.Range("B19").Select()
